$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$src = $ws.Range("Q4")
$dst = $ws.Range("R4")
$dst.NumberFormat = $src.NumberFormat
$dst.Font.Name = $src.Font.Name
$dst.Font.Size = $src.Font.Size
$dst.Font.Bold = $src.Font.Bold
$dst.HorizontalAlignment = $src.HorizontalAlignment
$dst.VerticalAlignment = $src.VerticalAlignment
$dst.Borders.Item(9).LineStyle = $src.Borders.Item(9).LineStyle
$dst.Value = 2021

Write-Output "done"
